$d = $word.ActiveDocument

# Edit 1: paragraph 24 - split "_Client_0 = " run into "_Client_0" / "_1" / " = "
$p24 = $d.Paragraphs(24).Range
$xml24 = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:p w14:paraId="3B3B6BC8" w14:textId="5538A7E7" w:rsidR="0074785B" w:rsidRPr="003167B2" w:rsidRDefault="0074785B" w:rsidP="003167B2"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="003167B2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>train</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="003167B2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>_Client_0</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>_1</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="003167B2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>file.loc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="003167B2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>[file[''</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="003167B2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>client_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="003167B2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>''] == ''train_Client_0'']</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:xmlData>'
$p24.InsertXML($xml24)

# Edit 3: paragraph 29 - add lastRenderedPageBreak before the first "file" run
$p29 = $d.Paragraphs(29).Range
$xml29 = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:p w14:paraId="6D2433FA" w14:textId="77777777" w:rsidR="0074785B" w:rsidRPr="0074785B" w:rsidRDefault="0074785B" w:rsidP="003167B2"><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="0074785B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>file</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="0074785B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>[''</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="0074785B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>counter_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="0074785B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">''] = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="0074785B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>encoder.fit_transform</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="0074785B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>(file[''</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="0074785B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>counter_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="0074785B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>''])</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:xmlData>'
$p29.InsertXML($xml29)

# Edit 2: paragraph 25 - split "(train_Client_0)" run, then append 4 new paragraphs after it
$p25 = $d.Paragraphs(25).Range
$xml25 = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:p w14:paraId="2ED3EA0F" w14:textId="71A2B58A" w:rsidR="0074785B" w:rsidRDefault="0074785B" w:rsidP="003167B2"><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="0074785B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>print</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="0074785B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>(train_Client_0</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>_1</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>#Deuxième méthode</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>train</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">_Client_0_2 = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>file.iloc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>[0 :35]</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>print</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>(train_Client_0</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>_</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>)</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:xmlData>'
$p25.InsertXML($xml25)
